# Update cryptocurrency price/volume data to latest scraped values
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "'63.313.81"
$ws.Range("E2").Value = "  -1.44%  "

# Row 3
$ws.Range("D3").Value = "'2.682.14"
$ws.Range("E3").Value = "  -3.09%  "

# Row 4
$ws.Range("E4").Value = "  +0.04%  "

# Row 5
$ws.Range("D5").Value = "'551.55"
$ws.Range("E5").Value = "  -4.33%  "

# Row 6
$ws.Range("D6").Value = "'158.13"
$ws.Range("E6").Value = "  -0.87%  "

# Row 7
$ws.Range("D7").Value = "'1.00"
$ws.Range("E7").Value = "  +0.38%  "

# Row 8
$ws.Range("E8").Value = "  -2.46%  "

# Row 9
$ws.Range("D9").Value = "'0.105"
$ws.Range("E9").Value = "  -4.47%  "

# Row 10
$ws.Range("E10").Value = "  -1.54%  "

# Row 11
$ws.Range("E11").Value = "  -4.52%  "

# Row 12
$ws.Range("D12").Value = "'5.28"
$ws.Range("E12").Value = "  -10.19%  "

# Row 13
$ws.Range("D13").Value = "'3.158.99"
$ws.Range("E13").Value = "  -2.86%  "

# Row 14
$ws.Range("D14").Value = "'26.33"
$ws.Range("E14").Value = "  -2.19%  "

# Row 15
$ws.Range("D15").Value = "'63.186.58"
$ws.Range("E15").Value = "  -1.00%  "

# Row 16
$ws.Range("D16").Value = "'0.0000145"
$ws.Range("E16").Value = "  -4.64%  "

# Row 17
$ws.Range("D17").Value = "'2.686.42"
$ws.Range("E17").Value = "  -3.03%  "

# Row 18
$ws.Range("D18").Value = "'11.98"
$ws.Range("E18").Value = "  -1.80%  "

# Row 19
$ws.Range("D19").Value = "'4.56"
$ws.Range("E19").Value = "  -5.75%  "

# Row 20
$ws.Range("D20").Value = "'346.33"
$ws.Range("E20").Value = "  -4.32%  "

# Row 21
$ws.Range("D21").Value = "'6.30"
$ws.Range("E21").Value = "  -5.16%  "

# Row 22
$ws.Range("E22").Value = "  -0.15%  "

# Row 23
$ws.Range("D23").Value = "'0.506"
$ws.Range("E23").Value = "  -4.23%  "

# Row 24
$ws.Range("D24").Value = "'63.81"
$ws.Range("E24").Value = "  -1.81%  "

# Row 25
$ws.Range("E25").Value = "  -1.37%  "

# Row 26
$ws.Range("E26").Value = "  +0.08%  "

# Row 27
$ws.Range("D27").Value = "'8.16"
$ws.Range("E27").Value = "  -4.75%  "

# Row 28
$ws.Range("D28").Value = "'0.0₃0856"
$ws.Range("E28").Value = "  -6.45%  "

# Row 29
$ws.Range("B29").Value = "PancakeSwap"
$ws.Range("C29").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D29").Value = "'1.93"
$ws.Range("E29").Value = "  -2.28%  "

# Row 30
$ws.Range("B30").Value = "Fetch.AI"
$ws.Range("C30").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D30").Value = "'1.34"
$ws.Range("E30").Value = "  -0.95%  "

# Row 31
$ws.Range("D31").Value = "'7.01"
$ws.Range("E31").Value = "  -5.22%  "

# Row 32
$ws.Range("D32").Value = "'166.21"
$ws.Range("E32").Value = "  -0.85%  "

# Row 33
$ws.Range("D33").Value = "'0.999"
$ws.Range("E33").Value = "  +0.07%  "

# Row 34
$ws.Range("D34").Value = "'4.81"
$ws.Range("E34").Value = "  -3.09%  "

# Row 35
$ws.Range("D35").Value = "'19.54"
$ws.Range("E35").Value = "  -3.34%  "

# Row 36
$ws.Range("D36").Value = "'1.43"
$ws.Range("E36").Value = "  -5.58%  "

# Row 37
$ws.Range("D37").Value = "'1.77"
$ws.Range("E37").Value = "  -2.51%  "

# Row 38
$ws.Range("D38").Value = "'340.57"
$ws.Range("E38").Value = "  -2.50%  "

# Row 39
$ws.Range("D39").Value = "'0.938"
$ws.Range("E39").Value = "  -6.72%  "

# Row 40
$ws.Range("D40").Value = "'6.13"
$ws.Range("E40").Value = "  -3.03%  "

# Row 41
$ws.Range("D41").Value = "'38.12"
$ws.Range("E41").Value = "  -2.79%  "

# Row 42
$ws.Range("E42").Value = "  -6.20%  "

# Row 43
$ws.Range("D43").Value = "'20.30"
$ws.Range("E43").Value = "  -5.77%  "

# Row 44
$ws.Range("D44").Value = "'20.75"
$ws.Range("E44").Value = "  -7.30%  "

# Row 45
$ws.Range("D45").Value = "'0.618"
$ws.Range("E45").Value = "  -2.28%  "

# Row 46
$ws.Range("D46").Value = "'0.0561"
$ws.Range("E46").Value = "  -5.47%  "

# Row 47
$ws.Range("E47").Value = "  +0.14%  "

# Row 48
$ws.Range("D48").Value = "'11.07"
$ws.Range("E48").Value = "  +0.38%  "

# Row 49
$ws.Range("B49").Value = "Stellar"
$ws.Range("C49").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D49").Value = "'0.0972"
$ws.Range("E49").Value = "  -3.86%  "

# Row 50
$ws.Range("B50").Value = "Aave"
$ws.Range("C50").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D50").Value = "'129.30"
$ws.Range("E50").Value = "  -5.99%  "

# Row 51
$ws.Range("B51").Value = "Maker"
$ws.Range("C51").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D51").Value = "'2.092.93"
$ws.Range("E51").Value = "  -1.92%  "
